# debug club name matching in round participation calc
#
# The "Manual Calcs" column (F) on the two results sheets was a scratch
# column used while debugging club-name matching; remove it so "Club Name"
# shifts into F and "Performance points..." shifts into G (and any
# overflow club-name notes shift left by one column as well).

$wb = $excel.ActiveWorkbook

$ws703 = $wb.Worksheets.Item("Sydney Round 1 70.3")
$ws703.Columns("F").Delete()

$wsSprint = $wb.Worksheets.Item("Sydney Round 1 Sprint")
$wsSprint.Columns("F").Delete()
